# Auto-generated edit script: update column F ('想去人数') values
# across the 展览, 本地生活, and 全部类型 sheets per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 300
$ws.Range("F6").Value = 391
$ws.Range("F7").Value = 846
$ws.Range("F8").Value = 39
$ws.Range("F9").Value = 500
$ws.Range("F11").Value = 292
$ws.Range("F12").Value = 125
$ws.Range("F14").Value = 227
$ws.Range("F15").Value = 29
$ws.Range("F17").Value = 6571
$ws.Range("F18").Value = 61
$ws.Range("F21").Value = 7513
$ws.Range("F24").Value = 3376
$ws.Range("F25").Value = 22
$ws.Range("F26").Value = 1166
$ws.Range("F27").Value = 879
$ws.Range("F28").Value = 4510
$ws.Range("F29").Value = 19
$ws.Range("F32").Value = 202
$ws.Range("F33").Value = 192
$ws.Range("F34").Value = 1603
$ws.Range("F39").Value = 1182
$ws.Range("F40").Value = 1699
$ws.Range("F41").Value = 2128

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1215

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1215
$ws.Range("F7").Value = 300
$ws.Range("F8").Value = 391
$ws.Range("F9").Value = 846
$ws.Range("F10").Value = 39
$ws.Range("F11").Value = 500
$ws.Range("F14").Value = 292
$ws.Range("F15").Value = 125
$ws.Range("F18").Value = 227
$ws.Range("F19").Value = 29
$ws.Range("F21").Value = 6571
$ws.Range("F22").Value = 61
$ws.Range("F25").Value = 7513
$ws.Range("F28").Value = 3376
$ws.Range("F29").Value = 22
$ws.Range("F30").Value = 1166
$ws.Range("F31").Value = 879
$ws.Range("F32").Value = 4510
$ws.Range("F33").Value = 19
$ws.Range("F37").Value = 202
$ws.Range("F38").Value = 192
$ws.Range("F39").Value = 1603
$ws.Range("F44").Value = 1182
$ws.Range("F45").Value = 1699
$ws.Range("F47").Value = 2128
